$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.172.43'
$ws.Range("E2").Value = '  -3.28%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.863.70'
$ws.Range("E3").Value = '  -3.90%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9999'
$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '234.09'
$ws.Range("E5").Value = '  -3.32%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9998'
$ws.Range("E6").Value = '  -0.03%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4668'
$ws.Range("E7").Value = '  -2.67%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2832'
$ws.Range("E8").Value = '  -2.69%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06554'
$ws.Range("E9").Value = '  -3.38%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.20'
$ws.Range("E10").Value = '  +0.23%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07811'
$ws.Range("E11").Value = '  -0.35%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '96.16'
$ws.Range("E12").Value = '  -7.77%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.857.14'
$ws.Range("E13").Value = '  -4.34%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.133'
$ws.Range("E14").Value = '  -3.14%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6702'
$ws.Range("E15").Value = '  -4.05%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '280.24'
$ws.Range("E16").Value = '  -5.41%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '30.193.13'
$ws.Range("E17").Value = '  -3.20%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.000'
$ws.Range("E18").Value = '  +0.02%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.454'
$ws.Range("E19").Value = '  -2.09%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.64'
$ws.Range("E20").Value = '  -2.92%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.099.31'
$ws.Range("E21").Value = '  -4.90%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.000007259'
$ws.Range("E22").Value = '  -4.61%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.9999'
$ws.Range("E23").Value = '  -0.04%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.143'
$ws.Range("E24").Value = '  -4.31%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.316'
$ws.Range("E25").Value = '  -2.61%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '165.40'
$ws.Range("E26").Value = '  -2.30%  '

$ws.Range("E27").Value = '  -4.55%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.906'
$ws.Range("E28").Value = '  -9.29%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.344'
$ws.Range("E29").Value = '  -3.41%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.09610'

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.403'
$ws.Range("E31").Value = '  -4.82%  '

$ws.Range("E32").Value = '  -4.34%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.120'
$ws.Range("E33").Value = '  -5.09%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04663'
$ws.Range("E34").Value = '  -3.70%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7013'
$ws.Range("E35").Value = '  -5.07%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.099'
$ws.Range("E36").Value = '  -3.18%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.711'
$ws.Range("E37").Value = '  -0.47%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01855'
$ws.Range("E38").Value = '  -5.38%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.273'
$ws.Range("E39").Value = '  -8.10%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.529'
$ws.Range("E40").Value = '  -3.87%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '72.48'
$ws.Range("E41").Value = '  -5.44%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.8526'
$ws.Range("E42").Value = '  -2.31%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.929'
$ws.Range("E43").Value = '  -5.25%  '

$ws.Range("B44").Value = 'TheSandbox'
$ws.Range("C44").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.4166'
$ws.Range("E44").Value = '  -4.76%  '

$ws.Range("B45").Value = 'PaxDollar'
$ws.Range("C45").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.9998'
$ws.Range("E45").Value = '  -0.01%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '103.19'
$ws.Range("E46").Value = '  -2.47%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '988.12'
$ws.Range("E47").Value = '  -3.94%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.138'
$ws.Range("E48").Value = '  -5.91%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.220'
$ws.Range("E49").Value = '  -0.77%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '34.12'
$ws.Range("E50").Value = '  -3.19%  '
